$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Mapping of (row, col) -> new text, using 1-based indices as in the
# Word object model. Only the data rows (1, 5, 9, 13, 17) contain text;
# the interleaved rows are blank spacer rows.
$updates = @(
    @{ Row = 1;  Col = 1; Text = "89÷9=9, 8" }
    @{ Row = 1;  Col = 2; Text = "60÷2=30, 0" }
    @{ Row = 1;  Col = 3; Text = "62÷7=8, 6" }
    @{ Row = 1;  Col = 4; Text = "60÷2=30, 0" }
    @{ Row = 1;  Col = 5; Text = "74÷8=9, 2" }

    @{ Row = 5;  Col = 1; Text = "56÷8=7, 0" }
    @{ Row = 5;  Col = 2; Text = "81÷8=10, 1" }
    @{ Row = 5;  Col = 3; Text = "87÷5=17, 2" }
    @{ Row = 5;  Col = 4; Text = "69÷2=34, 1" }
    @{ Row = 5;  Col = 5; Text = "37÷8=4, 5" }

    @{ Row = 9;  Col = 1; Text = "38÷3=12, 2" }
    @{ Row = 9;  Col = 2; Text = "15÷3=5, 0" }
    @{ Row = 9;  Col = 3; Text = "17÷9=1, 8" }
    @{ Row = 9;  Col = 4; Text = "88÷8=11, 0" }
    @{ Row = 9;  Col = 5; Text = "23÷7=3, 2" }

    @{ Row = 13; Col = 1; Text = "65÷5=13, 0" }
    @{ Row = 13; Col = 2; Text = "48÷8=6, 0" }
    @{ Row = 13; Col = 3; Text = "31÷7=4, 3" }
    @{ Row = 13; Col = 4; Text = "22÷8=2, 6" }
    @{ Row = 13; Col = 5; Text = "65÷3=21, 2" }

    @{ Row = 17; Col = 1; Text = "17÷2=8, 1" }
    @{ Row = 17; Col = 2; Text = "66÷8=8, 2" }
    @{ Row = 17; Col = 3; Text = "74÷6=12, 2" }
    @{ Row = 17; Col = 4; Text = "48÷7=6, 6" }
    @{ Row = 17; Col = 5; Text = "38÷7=5, 3" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
